$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46, column B ("politeness_score") was stored as text "3"; it should
# become a real number 3.
$ws.Range("B46").Value = 3

# Append a new annotation row (row 47) for Ruilin.
$ws.Range("A47").Value = "Ruilin"

# politeness_score "1" must stay textual (matches the source data, which
# keeps this column as text) rather than being auto-coerced to a number.
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "1"

$ws.Range("C47").Value = "rude and misleading, constant series of willful misinterpretations and falsehoods"
$ws.Range("D47").Value = "FBK"
$ws.Range("E47").Value = "OTH"
$ws.Range("F47").Value = "f6e31c12-680e-4edf-b736-d4a426f6f32f"
$ws.Range("G47").Value = "Syg-YfWCW_annotated.xlsx"
$ws.Range("H47").Value = "It is incredible that the commenter continues to be so rude and misleading (should OpenReview have a moderating system?), and continues to frame this interaction as an attempt to convince *them* rather than to correct the constant series of willful misinterpretations and falsehoods that they manage to state about our work in every single interaction, in the hope that they do not mislead others."
